$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Fusível interno")
$ws2 = $wb.Worksheets.Item("Aprovados")
$ws3 = $wb.Worksheets.Item("Tabela")

# --- Sheet1: Fusivel interno ---
$ws1.Range("B2").Value = 25.12562814070353
$ws1.Range("B3").Value = 8.372839021946804
$ws1.Range("B4").Value = 200.0000000000001
$ws1.Range("B6").Value = 11
$ws1.Range("B8").Value = 0.5305164769729848

# --- Sheet2: Aprovados ---
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 0.31
$ws2.Range("D2").Value = 38
$ws2.Range("E2").Value = 24.5
$ws2.Range("F2").Value = 69
$ws2.Range("G2").Value = 35.88039215686275
$ws2.Range("H2").Value = 0.7629239313978801
$ws2.Range("I2").Value = 1.183310587474263
$ws2.Range("J2").Value = 0.8400962150695708
$ws2.Range("K2").Value = 0.04878205554956144
$ws2.Range("L2").Value = 0.7629239313978801
$ws2.Range("M2").Value = 1.397935380692977
$ws2.Range("N2").Value = 35.88039215686275
$ws2.Range("O2").Value = 10.72470940855171

# --- Sheet3: Tabela ---
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = 1
$ws3.Range("D2").Value = 2.284444180880472
$ws3.Range("E2").Value = 0.09324261962777437
$ws3.Range("F2").Value = 1.397935380692977
$ws3.Range("G2").Value = 0.2688307888222627
$ws3.Range("B3").Value = 1.006578947368421
$ws3.Range("C3").Value = 1.085526315789474
$ws3.Range("D3").Value = 2.47982427529788
$ws3.Range("E3").Value = 0.1012173173590972
$ws3.Range("F3").Value = 1.482553309618666
$ws3.Range("G3").Value = 0.3167814328984634
$ws3.Range("B4").Value = 1.014388489208633
$ws3.Range("C4").Value = 1.18705035971223
$ws3.Range("D4").Value = 2.711750286656675
$ws3.Range("E4").Value = 0.1106836851696602
$ws3.Range("F4").Value = 1.575851798120855
$ws3.Range("G4").Value = 0.3788063881624192
$ws3.Range("B5").Value = 1.023809523809524
$ws3.Range("C5").Value = 1.309523809523809
$ws3.Range("D5").Value = 2.991534046391093
$ws3.Range("E5").Value = 0.1221034304649426
$ws3.Range("F5").Value = 1.678077193093751
$ws3.Range("G5").Value = 0.4610051792445262
$ws3.Range("B6").Value = 1.035398230088496
$ws3.Range("C6").Value = 1.460176991150443
$ws3.Range("D6").Value = 3.335692830489185
$ws3.Range("E6").Value = 0.1361507277750688
$ws3.Range("F6").Value = 1.788336944445124
$ws3.Range("G6").Value = 0.5731786534330099
$ws3.Range("B7").Value = 1.05
$ws3.Range("C7").Value = 1.65
$ws3.Range("D7").Value = 3.769332898452778
$ws3.Range("E7").Value = 0.1538503223858277
$ws3.Range("F7").Value = 1.902939536968315
$ws3.Range("G7").Value = 0.7318918225686103
$ws3.Range("B8").Value = 1.068965517241379
$ws3.Range("C8").Value = 1.896551724137931
$ws3.Range("D8").Value = 4.332566549945723
$ws3.Range("E8").Value = 0.1768394510181928
$ws3.Range("F8").Value = 2.01129822905886
$ws3.Range("G8").Value = 0.9669597338731801
$ws3.Range("B9").Value = 1.094594594594595
$ws3.Range("C9").Value = 2.22972972972973
$ws3.Range("D9").Value = 5.093693106017269
$ws3.Range("E9").Value = 0.2079058410619293
$ws3.Range("F9").Value = 2.085032363369227
$ws3.Range("G9").Value = 1.336544599285263
$ws3.Range("B10").Value = 1.131147540983607
$ws3.Range("C10").Value = 2.704918032786887
$ws3.Range("D10").Value = 2.284444180880472
$ws3.Range("E10").Value = 0.09324261962777437
$ws3.Range("F10").Value = 1.397935380692977
$ws3.Range("G10").Value = 1.344153944111314
$ws3.Range("B11").Value = 1.1875
$ws3.Range("C11").Value = 3.4375
$ws3.Range("D11").Value = 2.284444180880472
$ws3.Range("E11").Value = 0.09324261962777437
$ws3.Range("F11").Value = 1.397935380692977
$ws3.Range("G11").Value = 2.688307888222627

# Remove row 12 (Grupos paralelos reduced from 12 to 11)
$ws3.Rows.Item(12).Delete()
